$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark pair that originally sits right
#    after the title paragraph's pPr (before its first run).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Helper: find `searchText` starting the search at character offset
# `afterPos`, returning the matched Range (Range.Start/.End give the hit).
# ---------------------------------------------------------------------------
function FindAfter($searchText, $afterPos) {
    $r = $d.Range($afterPos, $d.Content.End)
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Not found: $searchText"
    }
    return $r
}

# ---------------------------------------------------------------------------
# 2. Locate the runs inside the "עמוד 3: ... " paragraph that need editing:
#       ", מתי הוקמה החברה, וחזון החברה."
#    becomes
#       ", שנת הקמת החברה וחזון החברה."
#    with a fresh "_GoBack" bookmark dropped right before the final
#    " וחזון החברה." run.
# ---------------------------------------------------------------------------
$rLabel = FindAfter "עמוד 3:" 0

$rMati   = FindAfter " מתי "        $rLabel.End
$rHukma  = FindAfter "הוקמה"        $rMati.End
$rHevra  = FindAfter " החברה"       $rHukma.End
$rComma  = FindAfter ","            $rHevra.End
$rVision = FindAfter " וחזון החברה." $rComma.End

# Drop temporary bookmarks at every run boundary in this stretch so that the
# engine's run-coalescing (which merges adjacent same-formatted runs it
# touches) cannot bleed past the spots that must stay distinct runs.
$d.Bookmarks.Add("_tmpB0", $d.Range($rMati.Start, $rMati.Start))     | Out-Null
$d.Bookmarks.Add("_tmpB1", $d.Range($rMati.End,   $rMati.End))       | Out-Null
$d.Bookmarks.Add("_tmpB2", $d.Range($rHukma.End,  $rHukma.End))      | Out-Null
$d.Bookmarks.Add("_tmpB3", $d.Range($rComma.Start, $rComma.Start))   | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($rComma.End, $rComma.End))      | Out-Null

# " מתי " -> " "
$r1 = $d.Range($d.Bookmarks("_tmpB0").Range.Start, $d.Bookmarks("_tmpB1").Range.Start)
$r1.Text = " "

# "הוקמה" -> "שנת הקמת"
$r2 = $d.Range($d.Bookmarks("_tmpB1").Range.Start, $d.Bookmarks("_tmpB2").Range.Start)
$r2.Text = "שנת הקמת"

# "," (the one right after " החברה") -> " "
$r3 = $d.Range($d.Bookmarks("_tmpB3").Range.Start, $d.Bookmarks("_GoBack").Range.Start)
$r3.Text = " "

# Clean up the scaffolding bookmarks - only "_GoBack" should remain.
$d.Bookmarks("_tmpB0").Delete()
$d.Bookmarks("_tmpB1").Delete()
$d.Bookmarks("_tmpB2").Delete()
$d.Bookmarks("_tmpB3").Delete()
